# Rename the original sheet to "ec_class" first.
$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item(1).Name = "ec_class"

# Duplicate it (inserted before) -> this will become "db query vs get".
$wb.Worksheets.Item("ec_class").Copy($wb.Worksheets.Item("ec_class"))
$wb.Worksheets.Item(1).Name = "db query vs get"

# ----- Fix up "ec_class" sheet first so its new shared string is created
#       before the ones on "db query vs get" (keeps shared-string order
#       matching insertion order: pre_class once, @base.all(), @base.get()). -----
$ec = $wb.Worksheets.Item("ec_class")
$ec.Range("E1").Value = "pt_element: construct pre_class once"
$ec.Columns.Item(6).ColumnWidth = 18.142857142857142
$ec.Range("E2").Select()

# ----- Build "db query vs get" content from the duplicated sheet -----
$new = $wb.Worksheets.Item("db query vs get")

# Drop column C (old "Total pt_element_accumulator" block) so column D's
# width definition (spacer) slides in to become column C, and old E/F data
# slides left into D/E.
$new.Columns.Item(3).Delete()
# Drop the now-shifted old column F (currently E) - not needed on this sheet.
$new.Columns.Item(5).Delete()

# Clear the leftover "Total time"/"Total pt_element_accumulator" header row
# (row 2 B/D) - the new sheet doesn't use that row for B/D.
$new.Range("B2").ClearContents()
$new.Range("D2").ClearContents()

# Header row 1: quote-prefixed literal labels ("@base.all()" must be
# created before "@base.get()" to land at shared-string index 8, then 9).
$new.Range("D1").Value = "'@base.all()"
$new.Range("B1").Value = "'@base.get()"

# Row 3-5: formulas (recalculated by the engine).
$new.Range("B3").Formula = "=AVERAGE(B6:B32)"
$new.Range("D3").Formula = "=AVERAGE(D6:D32)"
$new.Range("B4").Formula = "=STDEVA(B6:B32)"
$new.Range("D4").Formula = "=STDEVA(D6:D32)"
$new.Range("B5").Formula = "=B4/B3"
$new.Range("D5").Formula = "=D4/D3"

# Rows 6-10: raw sample data.
$new.Range("B6").Value = 0.275005
$new.Range("B7").Value = 0.23533599999999999
$new.Range("B8").Value = 0.224244999999999
$new.Range("B9").Value = 0.23674799999999999
$new.Range("B10").Value = 0.25717099999999998

$new.Range("D6").Value = 0.32842199999999999
$new.Range("D7").Value = 0.37503599999999998
$new.Range("D8").Value = 0.27442299999999997
$new.Range("D9").Value = 0.35225099999999898
$new.Range("D10").Value = 0.33001699999999901

# Leave "db query vs get" as the selected/active tab, cursor on D10.
$new.Range("D10").Select()
$new.Activate()
